$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 551.5862
$ws.Range("J17").Value = 542.75
$ws.Range("L17").Value = 1628.25
$ws.Range("N17").Value = -1964.25
$ws.Range("H19").Value = 1197.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1197.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1197.5
$ws.Range("N19").Value = -1547.5
$ws.Range("H28").Value = 843.8333
$ws.Range("I28").Value = 698.3333
$ws.Range("J28").Value = 989.3333
$ws.Range("K28").Value = 698.3333
$ws.Range("L28").Value = 989.3333
$ws.Range("M28").Value = -213.3333
$ws.Range("N28").Value = -1959.3333
$ws.Range("H33").Value = 95.71429000000001
$ws.Range("I33").Value = 95.71429000000001
$ws.Range("K33").Value = 95.71429000000001
$ws.Range("M33").Value = 133.28571
$ws.Range("H38").Value = 7063.625
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 7063.625
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 21190.875
$ws.Range("N38").Value = -21934.875
$ws.Range("H62").Value = 5259.875
$ws.Range("I62").Value = 4428.4
$ws.Range("J62").Value = 6645.6665
$ws.Range("K62").Value = 4428.4
$ws.Range("L62").Value = 6645.6665
$ws.Range("M62").Value = -3804.4
$ws.Range("N62").Value = -7893.6665
$ws.Range("H65").Value = 5259.875
$ws.Range("I65").Value = 4428.4
$ws.Range("J65").Value = 6645.6665
$ws.Range("K65").Value = 22142
$ws.Range("L65").Value = 33228.3325
$ws.Range("M65").Value = -19022
$ws.Range("N65").Value = -39468.3325
$ws.Range("H76").Value = 2999.5
$ws.Range("I76").Value = 2999.5
$ws.Range("K76").Value = 2999.5
$ws.Range("M76").Value = -2684.5
$ws.Range("H79").Value = 2999.5
$ws.Range("I79").Value = 2999.5
$ws.Range("K79").Value = 2999.5
$ws.Range("M79").Value = -1907.5
$ws.Range("H88").Value = 5265469.5
$ws.Range("I88").Value = 12501605
$ws.Range("K88").Value = 12501605
$ws.Range("M88").Value = -12501199
$ws.Range("H91").Value = 5265469.5
$ws.Range("I91").Value = 12501605
$ws.Range("K91").Value = 12501605
$ws.Range("M91").Value = -12500201
$ws.Range("H111").Value = 5321.2856
$ws.Range("I111").Value = 8487.5
$ws.Range("J111").Value = 1099.6666
$ws.Range("K111").Value = 25462.5
$ws.Range("L111").Value = 3298.9998
$ws.Range("M111").Value = -22395.5
$ws.Range("N111").Value = -9432.9998
$ws.Range("H116").Value = 10452.429
$ws.Range("J116").Value = 12335.8
$ws.Range("L116").Value = 12335.8
$ws.Range("N116").Value = -19219.8
$ws.Range("M19").ClearContents()
$ws.Range("M38").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2345
$ws.Range("I102").Value = 1720.4615
$ws.Range("K102").Value = 1720.4615
$ws.Range("M102").Value = -98.46149999999989

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17244656
$ws.Range("I20").Value = 33338104
$ws.Range("K20").Value = 33338104
$ws.Range("M20").Value = -33337857
$ws.Range("H134").Value = 2798.125
$ws.Range("I134").Value = 1812.8334
$ws.Range("K134").Value = 5438.5002
$ws.Range("M134").Value = -2903.5002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6581903.5
$ws.Range("I31").Value = 2739.3572
$ws.Range("J31").Value = 25003562
$ws.Range("K31").Value = 2739.3572
$ws.Range("L31").Value = 25003562
$ws.Range("M31").Value = -2444.3572
$ws.Range("N31").Value = -25004152
$ws.Range("H34").Value = 6581903.5
$ws.Range("I34").Value = 2739.3572
$ws.Range("J34").Value = 25003562
$ws.Range("K34").Value = 2739.3572
$ws.Range("L34").Value = 25003562
$ws.Range("M34").Value = -2537.3572
$ws.Range("N34").Value = -25003966
$ws.Range("H86").Value = 5352.8
$ws.Range("I86").Value = 5483.385
$ws.Range("K86").Value = 5483.385
$ws.Range("M86").Value = -4360.385
$ws.Range("H89").Value = 5352.8
$ws.Range("I89").Value = 5483.385
$ws.Range("K89").Value = 27416.925
$ws.Range("M89").Value = -21800.925
$ws.Range("H122").Value = 2269.1482
$ws.Range("I122").Value = 2050.762
$ws.Range("J122").Value = 3033.5
$ws.Range("K122").Value = 6152.286
$ws.Range("L122").Value = 9100.5
$ws.Range("M122").Value = -3702.286
$ws.Range("N122").Value = -14000.5
$ws.Range("H134").Value = 4495.143
$ws.Range("I134").Value = 4627.7915
$ws.Range("K134").Value = 13883.3745
$ws.Range("M134").Value = -11348.3745

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 664.8333
$ws.Range("I5").Value = 664.8333
$ws.Range("K5").Value = 1994.4999
$ws.Range("M5").Value = -1882.4999
$ws.Range("H121").Value = 5321525.5
$ws.Range("J121").Value = 110458.7
$ws.Range("L121").Value = 331376.1
$ws.Range("N121").Value = -333996.1
$ws.Range("H129").Value = 80503.42999999999
$ws.Range("I129").Value = 841.6667
$ws.Range("K129").Value = 2525.0001
$ws.Range("M129").Value = 2474.9999
$ws.Range("H135").Value = 664.8333
$ws.Range("I135").Value = 664.8333
$ws.Range("K135").Value = 5983.4997
$ws.Range("M135").Value = -3448.4997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1646.0667
$ws.Range("I107").Value = 428.66666
$ws.Range("K107").Value = 428.66666
$ws.Range("M107").Value = 1491.33334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 24749.875
$ws.Range("I40").Value = 24749.875
$ws.Range("K40").Value = 24749.875
$ws.Range("M40").Value = -24613.875
$ws.Range("H82").Value = 559.3158
$ws.Range("I82").Value = 525.2
$ws.Range("K82").Value = 525.2
$ws.Range("M82").Value = -164.2
$ws.Range("H85").Value = 559.3158
$ws.Range("I85").Value = 525.2
$ws.Range("K85").Value = 525.2
$ws.Range("M85").Value = 722.8
$ws.Range("H93").Value = 2674.96
$ws.Range("I93").Value = 2472.7368
$ws.Range("K93").Value = 2472.7368
$ws.Range("M93").Value = -1224.7368
$ws.Range("H136").Value = 5680.4
$ws.Range("I136").Value = 5850.0835
$ws.Range("J136").Value = 5001.6665
$ws.Range("K136").Value = 17550.2505
$ws.Range("L136").Value = 15004.9995
$ws.Range("M136").Value = -15000.2505
$ws.Range("N136").Value = -20104.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 42947.5
$ws.Range("J54").Value = 42947.5
$ws.Range("L54").Value = 42947.5
$ws.Range("N54").Value = -43987.5
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 602.5
$ws.Range("J107").Value = 594
$ws.Range("K107").Value = 1807.5
$ws.Range("L107").Value = 1782
$ws.Range("M107").Value = 112.5
$ws.Range("N107").Value = -5622
$ws.Range("H132").Value = 3965.4546
$ws.Range("I132").Value = 4062.1
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 12186.3
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -9656.299999999999
$ws.Range("N132").Value = -14057
$ws.Range("H136").Value = 4054.2727
$ws.Range("I136").Value = 4259.7
$ws.Range("K136").Value = 12779.1
$ws.Range("M136").Value = -10229.1
